$d = $word.ActiveDocument

# 1) Remove "grey " from the sentence so the wording reads
#    "...while, or silver." (the surviving words become "while, " / "or " / "silver.")
$null = $d.Content.Find.Execute(
    "while, grey or silver.", $true, $false, $false, $false, $false,
    $true, 1, $false, "while, or silver.", 2)

# 2) Locate the paragraph's index again (1-based, like the Word object model).
$targetIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Dealers should prefer to buy cars with low age with color black, while, or silver.*") {
        $targetIdx = $i
    }
}
$target = $d.Paragraphs($targetIdx)

# 3) Split the single run into three runs so that
#    "Dealers should prefer ... while, " / "or " / "silver." each get their
#    own <w:r>. Toggling a character attribute on the middle word forces the
#    COM host to materialize run boundaries around it; flipping Bold back
#    off leaves the formatting identical to its neighbours while keeping
#    the split.
$prefix = "Dealers should prefer to buy cars with low age with color black, while, "
$midStart = $target.Range.Start + $prefix.Length
$midEnd = $midStart + "or ".Length
$midRange = $d.Range($midStart, $midEnd)
$midRange.Bold = 1
$midRange.Bold = 0

# 4) Add a new bulleted paragraph right after it with the new recommendation.
$target = $d.Paragraphs($targetIdx)
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($targetIdx + 1)
$newPara.Range.Text = "The top priced cars are tesla, Ferrari, ram and Porsche. "
